# Update GDP per Capita values (Montenegro) and extend data through 2016.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New GDP per Capita values for years 1952..2016 (index 0 -> 1952, etc.)
$gdpValues = @(
    "1784",
    "1570",
    "1581",
    "2099",
    "1906",
    "2125",
    "2072",
    "2200",
    "2439",
    "2922",
    "2847",
    "3217",
    "3732",
    "3802",
    "3948",
    "3995",
    "4160",
    "4546",
    "4857",
    "4994",
    "5206",
    "5137",
    "5394",
    "5512",
    "5740",
    "6298",
    "6583",
    "6523",
    "7992",
    "7635",
    "7790",
    "7610",
    "7828",
    "7804",
    "8021",
    "7592",
    "7299",
    "7278",
    "6605",
    "6068.21515512111",
    "4605.29885365734",
    "2981.80459571749",
    "3073.92072623232",
    "3571.07237599629",
    "4663.16006868561",
    "5121.72820568113",
    "5522.10439400571",
    "5142.72906840197",
    "6046.14222843826",
    "6429.56114185043",
    "6897.7249603766",
    "7445.05011825545",
    "8193.74969155216",
    "9002.21820467631",
    "10307.6469729876",
    "12027.4293764959",
    "13550.3563815497",
    "13459.4991586417",
    "14509.7655644148",
    "15747",
    "15421",
    "16064",
    "16437",
    "17249",
    "17741"
)

$lastExistingRow = 60      # row holding year 2010 before this edit
$firstExistingRow = 2      # row holding year 1952

# The "GDP per Capita" figures are stored as text (not numbers) in column E,
# so a leading apostrophe is used to force text storage instead of letting
# Excel auto-convert the numeric-looking strings into numbers.
for ($row = $firstExistingRow; $row -le $lastExistingRow; $row++) {
    $value = $gdpValues[$row - $firstExistingRow]
    $ws.Cells.Item($row, 5).Value = "'" + $value
}

# Append new rows for years 2011..2016
$newYears = 2011..2016
$startRow = $lastExistingRow + 1
for ($i = 0; $i -lt $newYears.Count; $i++) {
    $row = $startRow + $i
    $year = $newYears[$i]
    $value = $gdpValues[$row - $firstExistingRow]
    $ws.Cells.Item($row, 1).Value = 499
    $ws.Cells.Item($row, 2).Value = "Montenegro"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year
    $ws.Cells.Item($row, 5).Value = "'" + $value
}
